# "832. Jezus overwinnaar (bonustrack)" - replace the old two-line
# "Oceans (Where Feet May Fail)" stub with the full Dutch lyric sheet.
$d = $word.ActiveDocument

# Paragraph 1 (was the "Oceans..." title) becomes the song title, Heading 1.
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "832. Jezus overwinnaar (bonustrack)"
$p1.Style = "Heading 1"

# A brand-new paragraph 2 is inserted for the subtitle, styled Heading 2.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Liedje van Stichting Opwekking"
$p2.Style = "Heading 2"

# Paragraph 3 reuses the old paragraph 2 ("Liedje van Hillsong United") and
# keeps the default (Normal) style, becoming the first lyric line.
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Waar U verschijnt wordt alles nieuw"

# The rest of the lyrics: one plain paragraph per line, each inserted right
# after the previous one (inherits the Normal style from paragraph 3, so no
# explicit style reset is required).
$lyrics = @(
    "Want U bevrijdt en geeft leven",
    "Elke storm verstilt",
    "Door de klank van Uw stem",
    "Alles buigt voor Koning Jezus",
    "U bent de held die voor ons strijdt",
    "U baant de weg van overwinning",
    "Elke vijand vlucht en ieder bolwerk valt neer",
    "Naam boven alle namen, Hoogste Heer",
    "Voor eeuwig is Uw heerschappij",
    "Uw troon staat onwankelbaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "De duisternis licht op door U",
    "De duivel is door U verslagen",
    "Dood waar is je macht?",
    "Waar is je prikkel gebleven?",
    "Jezus leeft en ik zal leven!",
    "De schepping knielt in diepst ontzag",
    "De hemel juicht voor onze Koning",
    "En de machten van de hel",
    "Weten wie er regeert",
    "Naam boven alle namen, Hoogste Heer",
    "Voor eeuwig is Uw heerschappij",
    "Uw troon staat onwankelbaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "Voor eeuwig is Uw heerschappij",
    "Uw troon staat onwankelbaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "Mmm-mmm-mmm",
    "Naam boven alle namen",
    "Naam boven alle namen",
    "Naam boven alle namen",
    "Naam boven alle namen",
    "Naam boven alle namen (mmm-mmm-mmm)",
    "Naam boven alle namen (namen)",
    "Naam boven alle namen (Jezus)",
    "Naam boven alle namen (Jezus)",
    "Naam boven alle namen",
    "Naam boven alle namen",
    "Naam boven alle namen",
    "Voor eeuwig is Uw heerschappij",
    "Uw troon staat onwankelbaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "Voor eeuwig is Uw heerschappij",
    "Uw troon staat onwankelbaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "Ongeevenaarde kracht",
    "Ligt in Uw grote naam",
    "Jezus Overwinnaar",
    "U bent Jezus Overwinnaar",
    "U bent Jezus Overwinnaar"
)

$idx = 3
foreach ($line in $lyrics) {
    $d.Paragraphs($idx).Range.InsertParagraphAfter()
    $idx = $idx + 1
    $d.Paragraphs($idx).Range.Text = $line
}

Write-Output ("paragraph count: " + $d.Paragraphs.Count)
